# Refresh cryptos list values (prices / 1h volume %) pulled from coinranking.com
# Row 47/48 coin entries (THORChain / RenderToken) also swap order this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.208.92'
$ws.Range("E2").Value = '  +0.53%  '

# Row 3
$ws.Range("D3").Value = '2.071.07'
$ws.Range("E3").Value = '  +0.58%  '

# Row 4
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").Value = '''251.38'
$ws.Range("E5").Value = '  +0.94%  '

# Row 6
$ws.Range("D6").Value = '''0.676'
$ws.Range("E6").Value = '  +4.08%  '

# Row 7
$ws.Range("D7").Value = '''62.83'
$ws.Range("E7").Value = '  +26.53%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$ws.Range("D9").Value = '''61.18'
$ws.Range("E9").Value = '  +1.94%  '

# Row 10
$ws.Range("D10").Value = '''0.385'
$ws.Range("E10").Value = '  +4.88%  '

# Row 11
$ws.Range("D11").Value = '''0.0806'
$ws.Range("E11").Value = '  +9.83%  '

# Row 12
$ws.Range("E12").Value = '  +2.95%  '

# Row 13
$ws.Range("D13").Value = '''16.00'
$ws.Range("E13").Value = '  +7.86%  '

# Row 14
$ws.Range("D14").Value = '2.369.01'
$ws.Range("E14").Value = '  +0.31%  '

# Row 15
$ws.Range("D15").Value = '''0.827'
$ws.Range("E15").Value = '  +0.74%  '

# Row 16
$ws.Range("D16").Value = '''5.46'
$ws.Range("E16").Value = '  +8.65%  '

# Row 17
$ws.Range("D17").Value = '2.068.11'
$ws.Range("E17").Value = '  -1.97%  '

# Row 18
$ws.Range("D18").Value = '37.145.60'
$ws.Range("E18").Value = '  +0.68%  '

# Row 19
$ws.Range("D19").Value = '''75.13'
$ws.Range("E19").Value = '  +4.95%  '

# Row 21
$ws.Range("D21").Value = '''15.08'
$ws.Range("E21").Value = '  +15.15%  '

# Row 22
$ws.Range("D22").Value = '''5.46'
$ws.Range("E22").Value = '  +6.62%  '

# Row 23
$ws.Range("D23").Value = '''240.11'
$ws.Range("E23").Value = '  +1.21%  '

# Row 24
$ws.Range("E24").Value = '  +0.14%  '

# Row 25
$ws.Range("E25").Value = '  -0.48%  '

# Row 26
$ws.Range("D26").Value = '''171.94'
$ws.Range("E26").Value = '  +2.27%  '

# Row 27
$ws.Range("D27").Value = '''9.26'
$ws.Range("E27").Value = '  +1.03%  '

# Row 28
$ws.Range("E28").Value = '  -1.15%  '

# Row 29
$ws.Range("D29").Value = '''2.04'
$ws.Range("E29").Value = '  +2.64%  '

# Row 30
$ws.Range("E30").Value = '  +3.94%  '

# Row 31
$ws.Range("D31").Value = '''1.10'
$ws.Range("E31").Value = '  +4.53%  '

# Row 32
$ws.Range("D32").Value = '''4.66'
$ws.Range("E32").Value = '  +5.09%  '

# Row 33
$ws.Range("D33").Value = '''0.0639'
$ws.Range("E33").Value = '  +6.64%  '

# Row 34
$ws.Range("E34").Value = '  +10.35%  '

# Row 35
$ws.Range("E35").Value = '  -0.83%  '

# Row 36
$ws.Range("E36").Value = '  +0.08%  '

# Row 37
$ws.Range("D37").Value = '''2.32'
$ws.Range("E37").Value = '  +2.68%  '

# Row 38
$ws.Range("E38").Value = '  -3.70%  '

# Row 39
$ws.Range("D39").Value = '''0.110'
$ws.Range("E39").Value = '  +26.86%  '

# Row 40
$ws.Range("D40").Value = '''1.36'
$ws.Range("E40").Value = '  +4.11%  '

# Row 41
$ws.Range("D41").Value = '''18.88'
$ws.Range("E41").Value = '  +8.73%  '

# Row 42
$ws.Range("E42").Value = '  +2.98%  '

# Row 43
$ws.Range("E43").Value = '  +2.31%  '

# Row 44
$ws.Range("D44").Value = '''98.40'
$ws.Range("E44").Value = '  +2.04%  '

# Row 45
$ws.Range("D45").Value = '''4.30'
$ws.Range("E45").Value = '  +26.60%  '

# Row 46
$ws.Range("D46").Value = '''2.79'
$ws.Range("E46").Value = '  +1.32%  '

# Row 47
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '''4.54'
$ws.Range("E47").Value = '  +16.71%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''2.54'
$ws.Range("E48").Value = '  +13.92%  '

# Row 49
$ws.Range("D49").Value = '1.309.07'
$ws.Range("E49").Value = '  +1.03%  '

# Row 50
$ws.Range("E50").Value = '  -0.30%  '

# Row 51
$ws.Range("D51").Value = '''6.93'
$ws.Range("E51").Value = '  +2.23%  '
